$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 159, shifting existing rows 159-198 down to 160-199
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new weekly price record
$ws.Cells.Item(159, 1).Value = 7
$ws.Cells.Item(159, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(159, 3).Value = "Ñuble"
$ws.Cells.Item(159, 4).Value = 44508
$ws.Cells.Item(159, 5).Value = 16
$ws.Cells.Item(159, 6).Value = 100114013
$ws.Cells.Item(159, 7).Value = "Zanahoria"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 100
$ws.Cells.Item(159, 11).Value = 7500
$ws.Cells.Item(159, 12).Value = 8000
$ws.Cells.Item(159, 13).Value = 7750
$ws.Cells.Item(159, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(159, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(159, 16).Value = 388
$ws.Cells.Item(159, 17).Value = 20
$ws.Cells.Item(159, 18).Value = "Hortaliza"
